$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 15: version 0.4.0 ---
# (cells are written in this particular order so new shared-string entries
# land at the same indexes the original author produced)
$ws.Range("D15").Value = "-Added graphic representation of population.`n-Fixed flexibility issue."
$ws.Range("A15").Value = "0.4.0"
$ws.Range("B15").Value = "AUTOMATA CELULAR - copia (21)`nNo EXE"
$ws.Range("C15").Value = $ws.Range("C14").Value()
$ws.Range("E15").Value = $ws.Range("E14").Value()
$ws.Range("F15").Value = $ws.Range("F14").Value()
$ws.Range("G15").Value = $ws.Range("G14").Value()

# --- Row 16: version 0.4.1 ---
$ws.Range("A16").Value = "0.4.1"
$ws.Range("D16").Value = "-Fixed group selection.`n-Fixed reciprocal association.`n-Graphic representation disabled."
$ws.Range("C16").Value = "-Make the code more readable using functions.`n-UI: Delete rows according to working functionality.`n*Graphic representation of F'.`n-Rework E calc.`n-Document every function.`n_OPTIONAL: Make it possible to reduce Niches on mutations.`n-When there is no data, the program crashes.`n-Check save data formatting.`n-Rework graphic representation to make an EXE."
$ws.Range("B16").Value = "AUTOMATA CELULAR - copia (25)"
$ws.Range("E16").Value = $ws.Range("E14").Value()
$ws.Range("F16").Value = $ws.Range("F14").Value()
$ws.Range("G16").Value = $ws.Range("G14").Value()

# Row heights are computed by Excel from the wrapped text of the tallest
# cell in each row; pin them explicitly to mirror Excel's real autofit result.
$ws.Rows.Item(15).RowHeight = 115.2
$ws.Rows.Item(16).RowHeight = 129.6

# Update the view: scrolled down and C16 selected, matching the saved file.
$ws.Range("C16").Select()
